$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the hours value for the last iteration row (B13: 2 -> 3); the
# SUM formula in B14 recalculates automatically.
$ws.Range("B13").Value = 3

# A13 was missing the thin border its neighbours (A5:A12) have; add it so
# the cell uses the same style and the now-unused style entry is dropped.
$ws.Range("A13").Borders.LineStyle = 1

# Update the view's selection/active cell as recorded for the sheet.
$ws.Range("E22").Select()
